$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet: "Burndown Chart" -> "Sprint Backlog" ---
$ws.Name = "Sprint Backlog"

# --- Fix a typo in the existing "Identify 3 Code Smells / By: Goncalo Rodrigues"
#     reviewer list: period -> comma before "Goncalo Gomes" (cell I6, style unchanged) ---
$ws.Range("I6").Value = "Identify 3 Code Smells `nBy: Gonçalo Rodrigues `nReviewers: Joana Cruz, Bárbara Correia, Gonçalo Gomes"

# --- Row 7: "Identify 3 Code Smells / By: Gonçalo Gomes" moves into the
#     "Reviewing" column (H/I) with its reviewer list ---
$ws.Range("H7").Value = "->"
$ws.Range("F7").Copy()
$ws.Range("H7").PasteSpecial(-4122)

$ws.Range("I7").Value = "Identify 3 Code Smells `nBy: Gonçalo Gomes`nReviewers: Bárbara Correia, Joana Cruz, Guilherme Santana"
$ws.Range("G7").Copy()
$ws.Range("I7").PasteSpecial(-4122)

# --- Row 13: "Identify 3 Pattern Designs / By: Gonçalo Gomes" moves into the
#     "Reviewing" column (H/I) with its reviewer list ---
$ws.Range("H13").Value = "->"
$ws.Range("F13").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$ws.Range("I13").Value = "Identify 3 Desing Patterns`nBy: Gonçalo Gomes`nReviewers: Bárbara Correia, Joana Cruz, Guilherme Santana"
$ws.Range("G13").Copy()
$ws.Range("I13").PasteSpecial(-4122)

# --- Row 14: "Identify 3 Pattern Designs / By: Guilherme Santana" moves into
#     the "Reviewing" column (H/I) with its reviewer list (H14 already had the
#     right style, only needs a value) ---
$ws.Range("H14").Value = "->"
$ws.Range("I14").Value = "Identify 3 Desing Patterns`nBy: Guilherme Santana`nReviewers: Bárbara Correia, Gonçalo Gomes, Gonçalo Rodrigues"
$ws.Range("G14").Copy()
$ws.Range("I14").PasteSpecial(-4122)

# --- Row 15: "Identify 3 Pattern Designs / By: Bárbara Correia" moves into
#     the "Reviewing" column (H/I) with its reviewer list ---
$ws.Range("H15").Value = "->"
$ws.Range("F15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("I15").Value = "Identify 3 Desing Patterns`nBy: Bárbara Correia`nReviewers: Gonçalo Rodrigues, Guilherme Santana, Joana Cruz"
$ws.Range("G15").Copy()
$ws.Range("I15").PasteSpecial(-4122)

# --- Pin every touched row to its final height (writing the new wrapped
#     Reviewing text can trigger an autofit reflow, so these are set last) ---
$ws.Rows.Item(7).RowHeight = 89.25
$ws.Rows.Item(8).RowHeight = 90.75
$ws.Rows.Item(9).RowHeight = 78
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(14).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 88.5
